# Author resized several of the input-table columns (A:E, G, I:J) so the
# Russian labels/values fit better, then left the selection on B2 having
# scrolled the sheet down to row 7 ("Входные данные").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Column width changes -------------------------------------------------
# (Target widths are expressed in the stored OOXML "character" width; the
# COM ColumnWidth property differs from the stored width by the standard
# 5/MaximumDigitWidth padding offset, which this runtime fixes at 5/6.)
$ws.Columns.Item(1).ColumnWidth  = 15.166666666666666   # -> 16
$ws.Columns.Item(2).ColumnWidth  = 8.451822916666666    # -> 9.28515625
$ws.Columns.Item(3).ColumnWidth  = 8.022135416666666    # -> 8.85546875
$ws.Columns.Item(4).ColumnWidth  = 8.307291666666666    # -> 9.140625
$ws.Columns.Item(5).ColumnWidth  = 6.736979166666667    # -> 7.5703125
$ws.Columns.Item(7).ColumnWidth  = 6.877604166666667    # -> 7.7109375
$ws.Columns.Item(9).ColumnWidth  = 7.307291666666667    # -> 8.140625
$ws.Columns.Item(10).ColumnWidth = 7.307291666666667    # -> 8.140625

# --- Selection / scroll position ------------------------------------------
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
